$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-27 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-28 Saturday", 2) | Out-Null
$d.Content.Find.Execute("47÷4=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=6, 6", 2) | Out-Null
$d.Content.Find.Execute("31÷7=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2) | Out-Null
$d.Content.Find.Execute("89÷5=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "39÷2=19, 1", 2) | Out-Null
$d.Content.Find.Execute("13÷4=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷2=38, 1", 2) | Out-Null
$d.Content.Find.Execute("91÷8=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "55÷2=27, 1", 2) | Out-Null
$d.Content.Find.Execute("35÷7=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=1, 2", 2) | Out-Null
$d.Content.Find.Execute("75÷2=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=18, 2", 2) | Out-Null
$d.Content.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "67÷2=33, 1", 2) | Out-Null
$d.Content.Find.Execute("63÷2=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷5=2, 4", 2) | Out-Null
$d.Content.Find.Execute("71÷4=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=10, 0", 2) | Out-Null
$d.Content.Find.Execute("43÷7=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷8=3, 2", 2) | Out-Null
$d.Content.Find.Execute("56÷9=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "44÷2=22, 0", 2) | Out-Null
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "25÷6=4, 1", 2) | Out-Null
$d.Content.Find.Execute("86÷2=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "66÷9=7, 3", 2) | Out-Null
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=14, 1", 2) | Out-Null
$d.Content.Find.Execute("40÷7=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "15÷5=3, 0", 2) | Out-Null
$d.Content.Find.Execute("35÷2=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "29÷5=5, 4", 2) | Out-Null
$d.Content.Find.Execute("43÷9=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=15, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷4=14, 1", 2) | Out-Null
$d.Content.Find.Execute("42÷4=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "92÷6=15, 2", 2) | Out-Null
$d.Content.Find.Execute("73÷7=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "62÷6=10, 2", 2) | Out-Null
$d.Content.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=3, 2", 2) | Out-Null
$d.Content.Find.Execute("53÷7=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "67÷8=8, 3", 2) | Out-Null
$d.Content.Find.Execute("60÷3=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷8=2, 3", 2) | Out-Null
